$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H38").Value = 2016772.9
$ws.Range("I38").Value = 3584363
$ws.Range("J38").Value = 1300
$ws.Range("K38").Value = 10753089
$ws.Range("L38").Value = 3900
$ws.Range("M38").Value = -10752717
$ws.Range("N38").Value = -4644
$ws.Range("H129").Value = 912.5
$ws.Range("J129").Value = 978.8421
$ws.Range("L129").Value = 2936.5263
$ws.Range("N129").Value = -12936.5263

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H5").Value = 209.77777
$ws.Range("I5").Value = 224.5
$ws.Range("J5").Value = 198
$ws.Range("K5").Value = 224.5
$ws.Range("L5").Value = 198
$ws.Range("M5").Value = -112.5
$ws.Range("N5").Value = -422
$ws.Range("H32").Value = 25779.857
$ws.Range("I32").Value = 4573.426
$ws.Range("J32").Value = 97351.56
$ws.Range("K32").Value = 4573.426
$ws.Range("L32").Value = 97351.56
$ws.Range("M32").Value = -4286.426
$ws.Range("N32").Value = -97925.56
$ws.Range("H37").Value = 12519.625
$ws.Range("I37").Value = 5001.25
$ws.Range("J37").Value = 20038
$ws.Range("K37").Value = 5001.25
$ws.Range("L37").Value = 20038
$ws.Range("M37").Value = -4728.25
$ws.Range("N37").Value = -20584
$ws.Range("H44").Value = 13922.5
$ws.Range("J44").Value = 13922.5
$ws.Range("L44").Value = 13922.5
$ws.Range("N44").Value = -14898.5
$ws.Range("H45").Value = 1624.1904
$ws.Range("I45").Value = 1635.9166
$ws.Range("J45").Value = 1608.5555
$ws.Range("K45").Value = 1635.9166
$ws.Range("L45").Value = 1608.5555
$ws.Range("M45").Value = -1258.9166
$ws.Range("N45").Value = -2362.5555
$ws.Range("H55").Value = 9084.286
$ws.Range("J55").Value = 9098.333000000001
$ws.Range("L55").Value = 9098.333000000001
$ws.Range("N55").Value = -9728.333000000001
$ws.Range("H63").Value = 2550
$ws.Range("I63").Value = 2275
$ws.Range("J63").Value = 3100
$ws.Range("K63").Value = 2275
$ws.Range("L63").Value = 3100
$ws.Range("M63").Value = -1589
$ws.Range("N63").Value = -4472
$ws.Range("H66").Value = 2550
$ws.Range("I66").Value = 2275
$ws.Range("J66").Value = 3100
$ws.Range("K66").Value = 11375
$ws.Range("L66").Value = 15500
$ws.Range("M66").Value = -7943
$ws.Range("N66").Value = -22364
$ws.Range("H80").Value = 14259.333
$ws.Range("J80").Value = 14259.333
$ws.Range("L80").Value = 14259.333
$ws.Range("N80").Value = -16255.333
$ws.Range("H83").Value = 14259.333
$ws.Range("J83").Value = 14259.333
$ws.Range("L83").Value = 42777.999
$ws.Range("N83").Value = -52761.999
$ws.Range("H88").Value = 1700.875
$ws.Range("I88").Value = 1433.3334
$ws.Range("J88").Value = 2503.5
$ws.Range("K88").Value = 1433.3334
$ws.Range("L88").Value = 2503.5
$ws.Range("M88").Value = -1027.3334
$ws.Range("N88").Value = -3315.5
$ws.Range("H91").Value = 1700.875
$ws.Range("I91").Value = 1433.3334
$ws.Range("J91").Value = 2503.5
$ws.Range("K91").Value = 1433.3334
$ws.Range("L91").Value = 2503.5
$ws.Range("M91").Value = -29.33339999999998
$ws.Range("N91").Value = -5311.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 209.77777
$ws.Range("I4").Value = 224.5
$ws.Range("J4").Value = 198
$ws.Range("K4").Value = 224.5
$ws.Range("L4").Value = 198
$ws.Range("M4").Value = -109.5
$ws.Range("N4").Value = -428
$ws.Range("H82").Value = 16293.333
$ws.Range("I82").Value = 2882.4285
$ws.Range("J82").Value = 35068.6
$ws.Range("K82").Value = 2882.4285
$ws.Range("L82").Value = 35068.6
$ws.Range("M82").Value = -2499.4285
$ws.Range("N82").Value = -35834.6
$ws.Range("H85").Value = 16293.333
$ws.Range("I85").Value = 2882.4285
$ws.Range("J85").Value = 35068.6
$ws.Range("K85").Value = 2882.4285
$ws.Range("L85").Value = 35068.6
$ws.Range("M85").Value = -1556.4285
$ws.Range("N85").Value = -37720.6
$ws.Range("H86").Value = 75753.336
$ws.Range("J86").Value = 2865.6667
$ws.Range("L86").Value = 2865.6667
$ws.Range("N86").Value = -5111.6667
$ws.Range("H89").Value = 75753.336
$ws.Range("J89").Value = 2865.6667
$ws.Range("L89").Value = 14328.3335
$ws.Range("N89").Value = -25560.3335

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H51").Value = 7071
$ws.Range("J51").Value = 7928.4287
$ws.Range("L51").Value = 7928.4287
$ws.Range("N51").Value = -9400.4287
$ws.Range("H59").Value = 24980
$ws.Range("J59").Value = 24980
$ws.Range("L59").Value = 24980
$ws.Range("N59").Value = -27270
$ws.Range("H60").Value = 19680
$ws.Range("J60").Value = 19680
$ws.Range("L60").Value = 19680
$ws.Range("N60").Value = -20702
$ws.Range("H61").Value = 7071
$ws.Range("J61").Value = 7928.4287
$ws.Range("L61").Value = 7928.4287
$ws.Range("N61").Value = -8624.4287
$ws.Range("H68").Value = 14957.167
$ws.Range("J68").Value = 14957.167
$ws.Range("L68").Value = 14957.167
$ws.Range("N68").Value = -16455.167
$ws.Range("H71").Value = 14957.167
$ws.Range("J71").Value = 14957.167
$ws.Range("L71").Value = 44871.501
$ws.Range("N71").Value = -52359.501
$ws.Range("H74").Value = 19480.908
$ws.Range("J74").Value = 19480.908
$ws.Range("L74").Value = 19480.908
$ws.Range("N74").Value = -21228.908
$ws.Range("H77").Value = 19480.908
$ws.Range("J77").Value = 19480.908
$ws.Range("L77").Value = 58442.724
$ws.Range("N77").Value = -67178.724

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 1242.2933
$ws.Range("I131").Value = 1060.5555
$ws.Range("J131").Value = 1267.0758
$ws.Range("K131").Value = 3181.6665
$ws.Range("L131").Value = 3801.2274
$ws.Range("M131").Value = 1858.3335
$ws.Range("N131").Value = -13881.2274

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 2794.5
$ws.Range("I80").Value = 3260
$ws.Range("J80").Value = 2018.6666
$ws.Range("K80").Value = 3260
$ws.Range("L80").Value = 2018.6666
$ws.Range("M80").Value = -2262
$ws.Range("N80").Value = -4014.6666
$ws.Range("H83").Value = 2794.5
$ws.Range("I83").Value = 3260
$ws.Range("J83").Value = 2018.6666
$ws.Range("K83").Value = 16300
$ws.Range("L83").Value = 10093.333
$ws.Range("M83").Value = -11308
$ws.Range("N83").Value = -20077.333

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 563377.75
$ws.Range("I46").Value = 480
$ws.Range("J46").Value = 1266999.9
$ws.Range("K46").Value = 480
$ws.Range("L46").Value = 1266999.9
$ws.Range("M46").Value = -292
$ws.Range("N46").Value = -1267375.9

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 537.2308
$ws.Range("I113").Value = 463.66666
$ws.Range("J113").Value = 600.2857
$ws.Range("K113").Value = 1390.99998
$ws.Range("L113").Value = 1800.8571
$ws.Range("M113").Value = 779.0000199999999
$ws.Range("N113").Value = -6140.8571
